$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repulled data updates to column F (dSF) for several rows.
$ws.Range("F6").Value = -1
$ws.Range("F12").Value = 2
$ws.Range("F29").Value = 1
$ws.Range("F38").Value = 2
$ws.Range("F40").Value = 0
$ws.Range("F43").Value = -2
$ws.Range("F44").Value = 2
$ws.Range("F46").Value = 0
$ws.Range("F53").Value = 1
$ws.Range("F56").Value = -5
$ws.Range("F58").Value = 0
$ws.Range("F59").Value = -2
